$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览 (Exhibition)
$ws2 = $wb.Worksheets.Item(2)   # 演出 (Performance)
$ws4 = $wb.Worksheets.Item(4)   # 全部类型 (All types)

# ---- Sheet 1 (展览): simple numeric F-column bumps ----
$ws1.Range("F5").Value = 8686
$ws1.Range("F7").Value = 11014
$ws1.Range("F8").Value = 92
$ws1.Range("F9").Value = 14
$ws1.Range("F15").Value = 296
$ws1.Range("F22").Value = 1866
$ws1.Range("F24").Value = 610
$ws1.Range("F25").Value = 351
$ws1.Range("F27").Value = 75
$ws1.Range("F28").Value = 592
$ws1.Range("F30").Value = 1255
$ws1.Range("F35").Value = 1422
$ws1.Range("F39").Value = 31
$ws1.Range("F42").Value = 368
$ws1.Range("F47").Value = 143
$ws1.Range("F48").Value = 129

# ---- Sheet 1 (展览): row 34 G-column ticket becomes unavailable (numeric -> text) ----
$ws1.Range("G34").Value = "不可售"

# ---- Sheet 2 (演出): simple numeric changes ----
$ws2.Range("G3").Value = 180
$ws2.Range("F14").Value = 25
$ws2.Range("F18").Value = 65
$ws2.Range("F23").Value = 392

# ---- Sheet 4 (全部类型): simple numeric changes ----
$ws4.Range("G7").Value = 180
$ws4.Range("F9").Value = 8686
$ws4.Range("F11").Value = 11014
$ws4.Range("F12").Value = 92
$ws4.Range("F15").Value = 296
$ws4.Range("F20").Value = 1866
$ws4.Range("F22").Value = 610
$ws4.Range("F23").Value = 351
$ws4.Range("F25").Value = 75
$ws4.Range("F27").Value = 592
$ws4.Range("F35").Value = 1422
$ws4.Range("F37").Value = 65
$ws4.Range("F41").Value = 368
$ws4.Range("F45").Value = 392
$ws4.Range("F48").Value = 143
$ws4.Range("F49").Value = 129

# ---- Sheet 4 (全部类型): rows 29-32 content overwrite ----
# A new event (周杰伦作品国风音乐会) was inserted at row 29, pushing the three
# events that used to occupy rows 29-31 down into rows 30-32 (row index / column A
# itself is untouched), and the event that used to be at row 32 (EXA同人展) drops out.

# Force the B-column (pure "YYYY-MM-DD" strings) to be stored as text so Excel
# does not auto-convert them into date serial numbers.
$bDates = $ws4.Range("B29:B32")
$bDates.NumberFormat = "@"
$ws4.Range("B29").Value = "2024-10-11"
$ws4.Range("C29").Value = "北京·官方唯一授权-周杰伦作品国风音乐会"
$ws4.Range("D29").Value = "西直门外大街135号  北展剧场"
$ws4.Range("E29").Value = "2024.10.11 19:30-10.11 21:00"
$ws4.Range("F29").Value = 20
$ws4.Range("G29").Value = 180
$ws4.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=88666"
$ws4.Range("I29").Value = "//i1.hdslb.com/bfs/openplatform/202407/2KgWinEn1720077808243.jpeg"

$ws4.Range("B30").Value = "2024-10-19"
$ws4.Range("C30").Value = "北京·ICOS国际动漫节×CGF中国游戏节04"
$ws4.Range("D30").Value = "石景山路68号 北京首钢会展中心"
$ws4.Range("E30").Value = "2024.10.19 09:00-10.20 17:00"
$ws4.Range("F30").Value = 1255
$ws4.Range("G30").Value = 80
$ws4.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=88085"
$ws4.Range("I30").Value = "//i2.hdslb.com/bfs/openplatform/202406/jQr9LeQO1719381394199.jpeg"

$ws4.Range("B31").Value = "2024-10-19"
$ws4.Range("C31").Value = "北京·可行中国动漫游戏节"
$ws4.Range("D31").Value = "焦化路甲18号 东进国际中心"
$ws4.Range("E31").Value = "2024.10.19 09:00-10.20 18:00"
$ws4.Range("F31").Value = 23
$ws4.Range("G31").Value = 85
$ws4.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=92495"
$ws4.Range("I31").Value = "//i1.hdslb.com/bfs/openplatform/202409/28QBTqAo1726293348310.jpeg"

$ws4.Range("B32").Value = "2024-10-25"
$ws4.Range("C32").Value = "北京·伦敦西区音乐剧明星演唱会-经典版"
$ws4.Range("D32").Value = "西直门外大街135号（北京展览馆内） 北京展览馆剧场"
$ws4.Range("E32").Value = "2024.10.25 19:30-10.26 21:30"
$ws4.Range("F32").Value = 11
$ws4.Range("G32").Value = 144
$ws4.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=89359"
$ws4.Range("I32").Value = "//i0.hdslb.com/bfs/openplatform/202407/PzPiEKUI1721114840552.jpeg"

# Strip the auto-applied "Text" number format back off the date cells so the
# cells carry no style attribute, matching every other date cell in the sheet.
$bDates.ClearFormats()
